$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the disbursement date (Q2) and the payment schedule label (K2)
$ws.Range("Q2").Value = "20/12/2021"
$ws.Range("K2").Value = "Cronograma Pagos"

# Clear the saved cell selection on the sheet view (select A1 instead of F10)
$ws.Range("A1").Select()
